$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the existing "2022-Q2"
#    sheet (so it inherits identical column layout / header styling), then
#    place it right after "总计" and trim it down to a single data row.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item(1))

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Drop the extra data rows (3-7) that were copied from 2022-Q2, keeping just
# the header row and a single data row.
$q4.Range("A3:H7").EntireRow.Delete()

# Fill in the one data row for 2022-Q4.
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'001044"
$q4.Range("B2").Style = "Normal"
$q4.Range("C2").Value = "'嘉实新消费股票"
$q4.Range("C2").Style = "Normal"
$q4.Range("D2").Value = "'10.68"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").Value = "'81.38"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").Value = "'4.47"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").Value = "'0.4774"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# 2. Update the "总计" (Total) summary sheet: insert the new 2022-Q4 row at
#    the top of the data and shift every other quarter down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$rows = @(
    @(0, "2022-Q4", 1, 0.48),
    @(1, "2022-Q2", 6, 0.12),
    @(2, "2022-Q1", 1, 0.08),
    @(3, "2021-Q4", 5, 0.7),
    @(4, "2021-Q2", 9, 1.77),
    @(5, "2021-Q1", 6, 1.31),
    @(6, "2020-Q4", 8, 1.06)
)

# Make sure row 8 (brand new) has the same formatting as the other data
# rows before we populate it.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $total.Cells.Item($r, 1).Value = $data[0]
    $total.Cells.Item($r, 2).Value = $data[1]
    $total.Cells.Item($r, 3).Value = $data[2]
    $total.Cells.Item($r, 4).Value = $data[3]
}

# ---------------------------------------------------------------------------
# 3. Restore "总计" as the active sheet/view.
# ---------------------------------------------------------------------------
$total.Activate()
